# Excel_sheet_for_DataDrivenTesting.xlsx — "Updated POM, clean up code"
#
# Semantic content changes applied by this edit:
#   1. The old generic "password_mismatch:The two password fields didn't
#      match." message (used verbatim in Register!D6:D11) is replaced by six
#      distinct, more specific validation messages.
#   2. Register!B5:C5 keeps its existing merge; nothing else changes there.
#   3. phythoncode!B4 gains a "Please enter valid data" helper message that
#      used to be blank.
#   4. The new Register!D6:D11 messages get a light-grey boxed look
#      (medium grey border all around) with vertical-center / wrap-text
#      alignment to make the longer text readable.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Register sheet — replace the single reused "password mismatch" message
#    with six specific validation messages.
# ---------------------------------------------------------------------------
$register = $wb.Worksheets.Item("Register")

$register.Range("D6").Value  = "characters other than Letters, digits and @/./+/-/_ are not allowed"
$register.Range("D7").Value  = "Passwords cannot be numeric"
$register.Range("D8").Value  = "password must contain at least 8 characters"
$register.Range("D9").Value  = "Passwords cannot be different"
$register.Range("D10").Value = "Passwords cannot be similar to full name or any additions to name"
$register.Range("D11").Value = "Common passwords cannot be used"

# Give the new, longer messages a readable boxed style: a medium light-grey
# border around the cell, or wrapped text for the longer ones.
foreach ($addr in @("D6","D10")) {
    $rng = $register.Range($addr)
    $rng.Borders.LineStyle = 1
    $rng.Borders.Weight = -4138
    $rng.Borders.Color = 13421772
    $rng.VerticalAlignment = -4108
}

foreach ($addr in @("D7","D8","D9","D11")) {
    $rng = $register.Range($addr)
    $rng.Borders.LineStyle = 1
    $rng.Borders.Weight = -4138
    $rng.Borders.Color = 13421772
    $rng.WrapText = $true
}

# ---------------------------------------------------------------------------
# 2) phythoncode sheet — add the missing helper message under the sample
#    code, which used to be an empty cell.
# ---------------------------------------------------------------------------
$pythoncode = $wb.Worksheets.Item("phythoncode")
$pythoncode.Range("B4").Value = "Please enter valid data"
